$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update percentage values (B2, B3): 3 -> 1.3
$ws.Range("B2").Value = 1.3
$ws.Range("B3").Value = 1.3

# Update contribution assessment ceilings (B4, B5)
$ws.Range("B4").Value = 90600
$ws.Range("B5").Value = 89400

# Update the entry date string (B6): 15.12.2023 -> 01.01.2024
$ws.Range("B6").Value = "01.01.2024"

# Update selection to match the saved view state
$ws.Range("B7").Select()
